# Update Sheet2 contents: change "ee" -> "Yahoo" in C4 and "AOT" -> "Academy" in C7.
# Also move the active selection to C4 (matching the edited cell) as recorded in the sheetView.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("C7").Value = "Academy"
$ws.Range("C4").Value = "Yahoo"

$ws.Activate()
$ws.Range("C4").Select()
